$wb = $excel.ActiveWorkbook

# --- Create the new first sheet "test-data.csv-2" by duplicating
#     "test-data.csv-1" (this keeps page setup / extLst and produces the
#     correct rId/sheetId/tabSelected rewiring that Excel performs on copy). ---
$src = $wb.Worksheets.Item("test-data.csv-1")
$src.Copy($src)
$new = $wb.Worksheets.Item(1)
$new.Name = "test-data.csv-2"

# Remove the copied content; we'll repopulate with the new column layout.
$new.Cells.Clear()

# --- Headers (row 1) ---
$new.Range("A1").Value = "Subject"
$new.Range("B1").Value = "Session"
$new.Range("C1").Value = "TrialName"
$new.Range("D1").Value = "cueSlide.RT"
$new.Range("E1").Value = "cueDur"
$new.Range("F1").Value = "TrialTypeBG"
$new.Range("G1").Value = "TrialTypeFG"
$new.Range("H1").Value = "targetSlide.RT"

# --- Row 2 ---
$new.Range("A2").Value = 1
$new.Range("B2").Value = 2
$new.Range("C2").Value = "practice"
$new.Range("D2").Value = 0
$new.Range("E2").Value = 1000
$new.Range("F2").Value = "Congruent"
$new.Range("G2").Value = "Congruent"
$new.Range("H2").Value = 705

# --- Row 3 ---
$new.Range("A3").Value = 82
$new.Range("B3").Value = 2
$new.Range("C3").Value = "test"
$new.Range("D3").Value = 0
$new.Range("E3").Value = 1000
$new.Range("F3").Value = "Congruent"
$new.Range("G3").Value = "Incongruent"
$new.Range("H3").Value = 392

# --- Column widths (best fit to content, matching the recorded layout) ---
$new.Range("A1").EntireColumn.ColumnWidth = 6.333333333333333
$new.Range("B1").EntireColumn.ColumnWidth = 6.5
$new.Range("C1").EntireColumn.ColumnWidth = 8.833333333333334
$new.Range("D1").EntireColumn.ColumnWidth = 9.666666666666666
$new.Range("E1").EntireColumn.ColumnWidth = 6.333333333333333
$new.Range("F1").EntireColumn.ColumnWidth = 10.333333333333334
$new.Range("G1").EntireColumn.ColumnWidth = 10.333333333333334
$new.Range("H1").EntireColumn.ColumnWidth = 11.833333333333334

# Restore the tab selection/active-cell state recorded for this sheet.
$new.Activate()
$new.Range("J31").Select()
